$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Jun_13" data column (C),
# shifting that column (and everything in it) to E.
$ws.Columns("C:D").Insert()

# Keep the inserted/shifted columns at the same custom width (8.0 chars)
# as the original data column.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14

# Row 1 holds the date headers. Shift the header text right by two slots
# and introduce the two new snapshot dates (in creation order Jun_15 then
# Jun_17, mirroring how the original workbook's string table grew).
$ws.Range("D1").Value = $ws.Range("B1").Value2
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Every other row gets blank ("UN") placeholders in the two new columns,
# matching the existing "no update" placeholder used elsewhere in the sheet.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}
